{"js": "// Replace two-digit multiplication expressions per the target diff.\n// Each \"before\" string is unique in the document body, so a plain\n// contains-search + full-match replace is unambiguous and order-independent.\nconst pairs = [\n  [\"14\u00d719=\", \"28\u00d725=\"],\n  [\"57\u00d788=\", \"23\u00d778=\"],\n  [\"86\u00d734=\", \"50\u00d721=\"],\n  [\"12\u00d750=\", \"89\u00d753=\"],\n  [\"73\u00d728=\", \"12\u00d783=\"],\n  [\"90\u00d750=\", \"13\u00d798=\"],\n  [\"67\u00d795=\", \"67\u00d770=\"],\n  [\"28\u00d736=\", \"22\u00d742=\"],\n  [\"69\u00d790=\", \"83\u00d749=\"],\n  [\"47\u00d733=\", \"79\u00d785=\"],\n  [\"74\u00d759=\", \"21\u00d793=\"],\n  [\"81\u00d727=\", \"58\u00d793=\"],\n  [\"21\u00d774=\", \"50\u00d713=\"],\n  [\"70\u00d782=\", \"24\u00d791=\"],\n  [\"92\u00d735=\", \"98\u00d738=\"],\n  [\"70\u00d796=\", \"47\u00d726=\"],\n  [\"66\u00d726=\", \"20\u00d767=\"],\n  [\"90\u00d733=\", \"62\u00d735=\"],\n  [\"22\u00d781=\", \"13\u00d768=\"],\n  [\"72\u00d790=\", \"89\u00d758=\"],\n  [\"52\u00d711=\", \"75\u00d745=\"],\n  [\"72\u00d731=\", \"75\u00d714=\"],\n  [\"38\u00d790=\", \"27\u00d744=\"],\n  [\"58\u00d776=\", \"63\u00d735=\"],\n  [\"44\u00d715=\", \"82\u00d715=\"],\n];\n\nconst body = context.document.body;\nfor (const [findText, replaceText] of pairs) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the two-digit-by-two-digit multiplication prompts in the table.\n# Each \"before\" expression is unique within the document body, so a\n# Find/Replace (wdReplaceAll) per pair is unambiguous and order-independent.\n$pairs = @(\n    @(\"14\u00d719=\", \"28\u00d725=\"),\n    @(\"57\u00d788=\", \"23\u00d778=\"),\n    @(\"86\u00d734=\", \"50\u00d721=\"),\n    @(\"12\u00d750=\", \"89\u00d753=\"),\n    @(\"73\u00d728=\", \"12\u00d783=\"),\n    @(\"90\u00d750=\", \"13\u00d798=\"),\n    @(\"67\u00d795=\", \"67\u00d770=\"),\n    @(\"28\u00d736=\", \"22\u00d742=\"),\n    @(\"69\u00d790=\", \"83\u00d749=\"),\n    @(\"47\u00d733=\", \"79\u00d785=\"),\n    @(\"74\u00d759=\", \"21\u00d793=\"),\n    @(\"81\u00d727=\", \"58\u00d793=\"),\n    @(\"21\u00d774=\", \"50\u00d713=\"),\n    @(\"70\u00d782=\", \"24\u00d791=\"),\n    @(\"92\u00d735=\", \"98\u00d738=\"),\n    @(\"70\u00d796=\", \"47\u00d726=\"),\n    @(\"66\u00d726=\", \"20\u00d767=\"),\n    @(\"90\u00d733=\", \"62\u00d735=\"),\n    @(\"22\u00d781=\", \"13\u00d768=\"),\n    @(\"72\u00d790=\", \"89\u00d758=\"),\n    @(\"52\u00d711=\", \"75\u00d745=\"),\n    @(\"72\u00d731=\", \"75\u00d714=\"),\n    @(\"38\u00d790=\", \"27\u00d744=\"),\n    @(\"58\u00d776=\", \"63\u00d735=\"),\n    @(\"44\u00d715=\", \"82\u00d715=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n\n    $found = $find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n"}
